$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header cell I1, re-using the same header style as H1 (bold, centered, bordered)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "公司"

# Update row 2 values. All of A2:I2 become the text "1" (kept as text, not numeric)
1..9 | ForEach-Object {
    $cell = $ws.Cells.Item(2, $_)
    $cell.NumberFormat = "@"
    $cell.Value = "1"
    $cell.Style = "Normal"
}
